# Auto update: 2025-12-03 08:54:10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (OKLO)
$ws.Range("D2").Value = 91.84
$ws.Range("E2").Value = 40.4
$ws.Range("F2").Value = 2.56
$ws.Range("K2").Value = 56.8
$ws.Range("N2").Value = 66.04328690552585

# Row 3 (NuScale / SMR)
$ws.Range("D3").Value = 18.91
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = -5.17
$ws.Range("K3").Value = 49.8
$ws.Range("N3").Value = 66.04328690552585
